$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.913.36"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "2.901.20"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'567.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.58%  "

$ws.Range("D6").Value = "'143.81"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.41%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.502"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").Value = "2.899.03"
$ws.Range("E9").Value = "  -1.31%  "

$ws.Range("E10").Value = "  -3.12%  "

$ws.Range("E11").Value = "  -2.78%  "

$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("E13").Value = "  -1.48%  "

$ws.Range("D14").Value = "'32.14"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").Value = "3.381.32"
$ws.Range("E16").Value = "  -1.30%  "

$ws.Range("D17").Value = "61.830.68"
$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.897.41"
$ws.Range("E18").Value = "  -1.52%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'6.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").Value = "'430.73"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.50%  "

$ws.Range("E21").Value = "  -4.32%  "

$ws.Range("E22").Value = "  -1.76%  "

$ws.Range("E23").Value = "  -1.27%  "

$ws.Range("D24").Value = "'78.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.56%  "

$ws.Range("D25").Value = "'11.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("D26").Value = "'10.07"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -9.08%  "

$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("E28").Value = "  -3.39%  "

$ws.Range("E29").Value = "  +8.27%  "

$ws.Range("D30").Value = "'7.02"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.75%  "

$ws.Range("E31").Value = "  -2.60%  "

$ws.Range("E32").Value = "  -6.10%  "

$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("E34").Value = "  -2.79%  "

$ws.Range("D35").Value = "'25.63"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.81%  "

$ws.Range("D36").Value = "'0.956"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.14%  "

$ws.Range("E37").Value = "  -3.02%  "

$ws.Range("D38").Value = "'48.83"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("D39").Value = "'2.84"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.46%  "

$ws.Range("E40").Value = "  -4.53%  "

$ws.Range("D41").Value = "'0.115"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.08%  "

$ws.Range("E42").Value = "  -2.63%  "

$ws.Range("D43").Value = "'40.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.65%  "

$ws.Range("E44").Value = "  -2.60%  "

$ws.Range("D45").Value = "2.702.84"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("E46").Value = "  -0.84%  "

$ws.Range("D47").Value = "'131.68"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.68%  "

$ws.Range("D48").Value = "'345.13"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.27%  "

$ws.Range("E50").Value = "  -1.05%  "

$ws.Range("D51").Value = "'21.56"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.14%  "

